# Timesheet "Signed Off" update.
# Fills in the Supervisor Name field, and the (previously blank) supervisor
# signature / date line at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor Name: -> Ankita Gangotra
$ws.Range("G6").Value = "Ankita Gangotra"

# Supervisor signature block (row 27): initials + sign-off date.
# D27 picks up the same short-date formatting already used for the
# employee's sign-off date in D25.
$ws.Range("A27").Value = "A.G"
$ws.Range("D27").Value = (Get-Date -Year 2014 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("D27").NumberFormat = "m/d/yy"

# Keep the same view the author left the sheet in: scrolled down so row 17
# is the first visible row, with D27:E27 (the date they just filled in)
# selected.
$ws.Range("D27:E27").Select()
$excel.ActiveWindow.ScrollRow = 17
